$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.991.61"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "3.850.66"
$ws.Range("E3").Value = "  +1.21%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'707.42"
$ws.Range("E5").Value = "  +1.22%  "
$ws.Range("D6").Value = "'172.71"
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").Value = "3.847.22"
$ws.Range("E7").Value = "  +1.18%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -0.58%  "
$ws.Range("E10").Value = "  +0.12%  "
$ws.Range("D11").Value = "'7.34"
$ws.Range("E11").Value = "  -0.71%  "
$ws.Range("E12").Value = "  -0.58%  "
$ws.Range("E13").Value = "  -0.41%  "
$ws.Range("D14").Value = "'36.69"
$ws.Range("E14").Value = "  +0.70%  "
$ws.Range("D15").Value = "4.498.22"
$ws.Range("D16").Value = "3.840.89"
$ws.Range("E16").Value = "  +1.03%  "
$ws.Range("D17").Value = "71.013.77"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("E19").Value = "  +0.73%  "
$ws.Range("D20").Value = "'17.36"
$ws.Range("E20").Value = "  -2.81%  "
$ws.Range("D21").Value = "'497.15"
$ws.Range("E21").Value = "  +2.96%  "
$ws.Range("E22").Value = "  -3.95%  "
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("D24").Value = "'85.55"
$ws.Range("E24").Value = "  +1.40%  "
$ws.Range("E25").Value = "  +2.03%  "
$ws.Range("D26").Value = "'10.67"
$ws.Range("E26").Value = "  +1.61%  "
$ws.Range("E27").Value = "  -1.91%  "
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").Value = "'3.21"
$ws.Range("E28").Value = "  +3.38%  "
$ws.Range("B29").Value = "Fetch.AI"
$ws.Range("C29").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D29").Value = "'2.11"
$ws.Range("E29").Value = "  -3.22%  "
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("E32").Value = "  -0.88%  "
$ws.Range("D33").Value = "'29.50"
$ws.Range("E33").Value = "  -0.07%  "
$ws.Range("E34").Value = "  -3.17%  "
$ws.Range("E35").Value = "  -0.63%  "
$ws.Range("D36").Value = "3.805.85"
$ws.Range("E36").Value = "  +1.41%  "
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("E39").Value = "  +7.20%  "
$ws.Range("D40").Value = "'6.05"
$ws.Range("E40").Value = "  +0.18%  "
$ws.Range("E41").Value = "  +5.80%  "
$ws.Range("D42").Value = "'3.34"
$ws.Range("E42").Value = "  -3.19%  "
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("D45").Value = "'0.000322"
$ws.Range("E45").Value = "  -1.85%  "
$ws.Range("D46").Value = "'163.51"
$ws.Range("E46").Value = "  +0.70%  "
$ws.Range("D47").Value = "'48.64"
$ws.Range("E47").Value = "  -0.41%  "
$ws.Range("D48").Value = "'1.39"
$ws.Range("E48").Value = "  +0.70%  "
$ws.Range("D49").Value = "'416.51"
$ws.Range("E49").Value = "  +1.62%  "
$ws.Range("E50").Value = "  -1.16%  "
$ws.Range("D51").Value = "'8.61"
$ws.Range("E51").Value = "  +0.33%  "
